$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab ("tc1") query text in B2: remove the trailing
# `Cohort` column (co.cohort_description) that is no longer part of the
# UBC01 cases query, and drop the now-trailing comma on the prior line.
$newCasesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN [''UBC01''] and diag.stage_of_disease in [''T3N0M0'', ''Not Applicable''] 
     OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'
$ws.Range("B2").Value = $newCasesQuery

# The row auto-sizes to the (now one line shorter) wrapped text.
$ws.Rows.Item(2).RowHeight = 304.5

# Restore/refresh the view state: scroll back up so B2 is the top-left
# visible cell, and select B2 (matching the saved workbook view).
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B2").Select()
